$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.835.88'
$ws.Range('E2').Value = '  +1.61%  '

$ws.Range('D3').Value = '1.766.10'
$ws.Range('E3').Value = '  +2.09%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.004'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '328.06'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.77%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.06%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4465'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -1.28%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3548'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.60%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07418'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.83%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.00'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.74%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.098'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.30%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.003'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.05%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.89'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.62%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.017'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.73%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.233'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.65%  '

$ws.Range('D16').Value = '1.772.09'
$ws.Range('E16').Value = '  +2.59%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '93.24'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.47%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001060'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.03%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06421'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.58%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.001'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.04%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.09'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.12%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.767'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.78%  '

$ws.Range('D23').Value = '27.898.24'
$ws.Range('E23').Value = '  +1.67%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.93%  '

$ws.Range('E25').Value = '  +2.56%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.65'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.17%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.37'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.56%  '

$ws.Range('D28').Value = '1.975.30'
$ws.Range('E28').Value = '  +2.60%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.157'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.95%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '124.95'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.41%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.104'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +6.10%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09202'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.08%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.634'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.31%  '

$ws.Range('E34').Value = '  +0.70%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '11.83'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.39%  '

$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06175'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.85%  '

$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02284'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +1.15%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2097'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.49%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.6316'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.64%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '4.954'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.23%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.186'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.16%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.394'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.89%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '7.869'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.18%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.23'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.94%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.747'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.43%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5866'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.55%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '122.50'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.63%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.955'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.07%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06903'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.10%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.133'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.28%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '72.90'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.98%  '
